$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.739.08'
$ws.Range("E2").Value = '  +4.30%  '
$ws.Range("D3").Value = '3.500.11'
$ws.Range("E3").Value = '  +2.63%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.82'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.42'
$ws.Range("E6").Value = '  +7.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D8").Value = '3.499.54'
$ws.Range("E8").Value = '  +2.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.579'
$ws.Range("E9").Value = '  +4.90%  '
$ws.Range("E10").Value = '  +0.68%  '
$ws.Range("E11").Value = '  +4.81%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.437'
$ws.Range("E12").Value = '  +3.09%  '
$ws.Range("D13").Value = '4.100.50'
$ws.Range("E13").Value = '  +2.52%  '
$ws.Range("E14").Value = '  -0.56%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.24'
$ws.Range("E15").Value = '  +4.29%  '
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000178'
$ws.Range("E16").Value = '  +3.20%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '66.720.18'
$ws.Range("E17").Value = '  +4.16%  '
$ws.Range("D18").Value = '3.505.26'
$ws.Range("E18").Value = '  +3.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.34'
$ws.Range("E19").Value = '  +4.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.96'
$ws.Range("E20").Value = '  +3.12%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '388.64'
$ws.Range("E21").Value = '  +3.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.98'
$ws.Range("E22").Value = '  +1.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.14'
$ws.Range("E23").Value = '  +2.32%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.532'
$ws.Range("E25").Value = '  +3.28%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000123'
$ws.Range("E26").Value = '  +5.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.36'
$ws.Range("E27").Value = '  +7.99%  '
$ws.Range("E28").Value = '  +2.65%  '
$ws.Range("E29").Value = '  +0.07%  '
$ws.Range("E30").Value = '  +5.82%  '
$ws.Range("E31").Value = '  +6.36%  '
$ws.Range("E32").Value = '  +2.75%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.56'
$ws.Range("E33").Value = '  +3.14%  '
$ws.Range("E34").Value = '  +5.45%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  +6.87%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.52'
$ws.Range("E37").Value = '  +2.24%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.884'
$ws.Range("E38").Value = '  +5.08%  '
$ws.Range("E39").Value = '  +4.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.67'
$ws.Range("E40").Value = '  +5.41%  '
$ws.Range("E41").Value = '  +2.61%  '
$ws.Range("E42").Value = '  +1.93%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.28'
$ws.Range("E43").Value = '  +2.16%  '
$ws.Range("D44").Value = '2.818.89'
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '43.03'
$ws.Range("E45").Value = '  -0.17%  '
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '26.54'
$ws.Range("E46").Value = '  +1.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '358.14'
$ws.Range("E47").Value = '  +6.34%  '
$ws.Range("E48").Value = '  +6.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0311'
$ws.Range("E49").Value = '  +2.71%  '
$ws.Range("E50").Value = '  +4.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.81'
$ws.Range("E51").Value = '  +13.43%  '
